$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.8203074518761176
$ws.Range("J2").Value = 0.8725723693674974
$ws.Range("M2").Value = 0.8092246666666667
$ws.Range("N2").Value = 2.427674
$ws.Range("O2").Value = 0.2312345204823145
$ws.Range("P2").Value = 0.2486288754648592
$ws.Range("Q2").Value = 0.04417476532866667
$ws.Range("R2").Value = 0.397572887958
$ws.Range("S2").Value = 0.1896834002826433
$ws.Range("T2").Value = 0.2169466869575486

$ws.Range("I3").Value = 0.8203074518761176
$ws.Range("J3").Value = 0.8725723693674974
$ws.Range("O3").Value = 0.0765802413191472
$ws.Range("P3").Value = 0.08234090326259635
$ws.Range("S3").Value = 0.06281934262056782
$ws.Range("T3").Value = 0.0718483970557036

$ws.Range("I4").Value = 0.8203074518761176
$ws.Range("J4").Value = 0.8725723693674974
$ws.Range("M4").Value = 0.7391253333333334
$ws.Range("N4").Value = 2.217376
$ws.Range("O4").Value = 0.2112037596847816
$ws.Range("P4").Value = 0.2270913233666331
$ws.Range("Q4").Value = 0.04034811282133333
$ws.Range("R4").Value = 0.363133015392
$ws.Range("S4").Value = 0.1732520179336791
$ws.Range("T4").Value = 0.1981536140928235

$ws.Range("I5").Value = 0.8203074518761176
$ws.Range("J5").Value = 0.8725723693674974
$ws.Range("M5").Value = 0.7345045
$ws.Range("N5").Value = 1.469009
$ws.Range("O5").Value = 0.2098833647140458
$ws.Range("P5").Value = 0.150447735452848
$ws.Range("Q5").Value = 0.0400958661505
$ws.Range("R5").Value = 0.240575196903
$ws.Range("S5").Value = 0.1721688880997648
$ws.Range("T5").Value = 0.131276536990066

$ws.Range("I6").Value = 0.8203074518761176
$ws.Range("J6").Value = 0.8725723693674974
$ws.Range("M6").Value = 0.9487306666666666
$ws.Range("N6").Value = 2.846192
$ws.Range("O6").Value = 0.2710981137997109
$ws.Range("P6").Value = 0.2914911624530634
$ws.Range("Q6").Value = 0.05179025836266666
$ws.Range("R6").Value = 0.466112325264
$ws.Range("S6").Value = 0.2223838029394626
$ws.Range("T6").Value = 0.2543471342713556

$ws.Range("G7").Value = 0.011958
$ws.Range("H7").Value = 0.023916
$ws.Range("I7").Value = 0.1796925481238824
$ws.Range("J7").Value = 0.1274276306325027
$ws.Range("M7").Value = 0.8092246666666667
$ws.Range("N7").Value = 2.427674
$ws.Range("O7").Value = 0.2312345204823145
$ws.Range("P7").Value = 0.2486288754648592
$ws.Range("Q7").Value = 0.009676708564
$ws.Range("R7").Value = 0.058060251384
$ws.Range("S7").Value = 0.04155112019967116
$ws.Range("T7").Value = 0.03168218850731058

$ws.Range("G8").Value = 0.011958
$ws.Range("H8").Value = 0.023916
$ws.Range("I8").Value = 0.1796925481238824
$ws.Range("J8").Value = 0.1274276306325027
$ws.Range("O8").Value = 0.0765802413191472
$ws.Range("P8").Value = 0.08234090326259635
$ws.Range("Q8").Value = 0.003204732042
$ws.Range("R8").Value = 0.019228392252
$ws.Range("S8").Value = 0.01376089869857938
$ws.Range("T8").Value = 0.01049250620689276

$ws.Range("G9").Value = 0.011958
$ws.Range("H9").Value = 0.023916
$ws.Range("I9").Value = 0.1796925481238824
$ws.Range("J9").Value = 0.1274276306325027
$ws.Range("M9").Value = 0.7391253333333334
$ws.Range("N9").Value = 2.217376
$ws.Range("O9").Value = 0.2112037596847816
$ws.Range("P9").Value = 0.2270913233666331
$ws.Range("Q9").Value = 0.008838460736000002
$ws.Range("R9").Value = 0.05303076441600001
$ws.Range("S9").Value = 0.03795174175110251
$ws.Range("T9").Value = 0.02893770927380954

$ws.Range("G10").Value = 0.011958
$ws.Range("H10").Value = 0.023916
$ws.Range("I10").Value = 0.1796925481238824
$ws.Range("J10").Value = 0.1274276306325027
$ws.Range("M10").Value = 0.7345045
$ws.Range("N10").Value = 1.469009
$ws.Range("O10").Value = 0.2098833647140458
$ws.Range("P10").Value = 0.150447735452848
$ws.Range("Q10").Value = 0.008783204810999999
$ws.Range("R10").Value = 0.035132819244
$ws.Range("S10").Value = 0.03771447661428103
$ws.Range("T10").Value = 0.01917119846278199

$ws.Range("G11").Value = 0.011958
$ws.Range("H11").Value = 0.023916
$ws.Range("I11").Value = 0.1796925481238824
$ws.Range("J11").Value = 0.1274276306325027
$ws.Range("M11").Value = 0.9487306666666666
$ws.Range("N11").Value = 2.846192
$ws.Range("O11").Value = 0.2710981137997109
$ws.Range("P11").Value = 0.2914911624530634
$ws.Range("Q11").Value = 0.011344921312
$ws.Range("R11").Value = 0.06806952787199999
$ws.Range("S11").Value = 0.0487143108602483
$ws.Range("T11").Value = 0.0371440281817078
